$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 23
$ws.Range("H2").Value = 23
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 9
$ws.Range("H3").Value = 9
$ws.Range("F4").Value = 4
$ws.Range("H4").Value = 4
$ws.Range("F9").Value = 6
$ws.Range("H9").Value = 6
$ws.Range("E15").Value = 141
$ws.Range("F15").Value = 72
$ws.Range("H15").Value = 72
$ws.Range("F17").Value = 39
$ws.Range("H17").Value = 39
$ws.Range("E18").Value = 91
$ws.Range("F18").Value = 34
$ws.Range("H18").Value = 34
$ws.Range("F24").Value = 12
$ws.Range("H24").Value = 12
$ws.Range("F29").Value = 9
$ws.Range("H29").Value = 9
$ws.Range("F36").Value = 30
$ws.Range("H36").Value = 30
$ws.Range("E40").Value = 14
$ws.Range("F41").Value = 14
$ws.Range("H41").Value = 14
$ws.Range("F42").Value = 12
$ws.Range("H42").Value = 12
$ws.Range("E45").Value = 21
$ws.Range("F45").Value = 11
$ws.Range("H45").Value = 11
$ws.Range("E48").Value = 21
$ws.Range("F49").Value = 30
$ws.Range("H49").Value = 30
$ws.Range("F50").Value = 5
$ws.Range("H50").Value = 5
$ws.Range("E53").Value = 5
$ws.Range("E62").Value = 37
$ws.Range("F62").Value = 9
$ws.Range("H62").Value = 9
$ws.Range("F63").Value = 6
$ws.Range("H63").Value = 6
$ws.Range("F66").Value = 16
$ws.Range("H66").Value = 16
$ws.Range("E68").Value = 12
$ws.Range("F68").Value = 6
$ws.Range("H68").Value = 6
$ws.Range("F69").Value = 8
$ws.Range("H69").Value = 8
$ws.Range("E70").Value = 33
$ws.Range("F70").Value = 14
$ws.Range("H70").Value = 14
$ws.Range("E71").Value = 25
$ws.Range("E76").Value = 39
$ws.Range("F76").Value = 13
$ws.Range("H76").Value = 13
$ws.Range("E77").Value = 45
$ws.Range("E78").Value = 38
$ws.Range("F78").Value = 13
$ws.Range("H78").Value = 13
$ws.Range("F88").Value = 9
$ws.Range("H88").Value = 9
